$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 273, shifting all existing data rows (old 273..364)
# down by one (new 274..365).
$ws.Rows.Item(273).Insert()

# Populate the newly inserted row 273 with the new weekly record.
# Columns H, I, N, Q keep the same values the (now shifted-down) old
# row 273 had; the rest of the fields carry the new observation.
$ws.Cells.Item(273,1).Value  = 9
$ws.Cells.Item(273,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(273,3).Value  = "Metropolitana"
$ws.Cells.Item(273,4).Value  = (Get-Date -Year 2022 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(273,5).Value  = 13
$ws.Cells.Item(273,6).Value  = 100112021
$ws.Cells.Item(273,7).Value  = "Ají"
$ws.Cells.Item(273,8).Value  = "Americana (o)"
$ws.Cells.Item(273,9).Value  = "Primera"
$ws.Cells.Item(273,10).Value = 25
$ws.Cells.Item(273,11).Value = 30000
$ws.Cells.Item(273,12).Value = 32000
$ws.Cells.Item(273,13).Value = 30960
$ws.Cells.Item(273,14).Value = '$/caja 25 kilos'
$ws.Cells.Item(273,15).Value = "Provincia de Limarí"
$ws.Cells.Item(273,16).Value = 1238
$ws.Cells.Item(273,17).Value = 25
$ws.Cells.Item(273,18).Value = "Hortaliza"
